# Apply the K-League 1 2023 dataset update:
#  - Four pairs of adjacent rows had their match data (columns F:V) swapped
#    back into the correct row order (the Indice / pais / torneio / temporada
#    / data_partida columns A:E were already correct and stay untouched).
#  - Two brand-new match rows (181 and 182) are appended at the end of the
#    sheet, copying the formatting of the last existing data row (180).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-RowFV($sheet, $row, $vals) {
    # Writes an 17-element array into columns F..V (6..22) of $row.
    for ($i = 0; $i -lt $vals.Length; $i++) {
        $sheet.Cells.Item($row, 6 + $i).Value = $vals[$i]
    }
}

function Set-RowFull($sheet, $row, $vals) {
    # Writes a 22-element array into columns A..V (1..22) of $row.
    for ($i = 0; $i -lt $vals.Length; $i++) {
        $sheet.Cells.Item($row, 1 + $i).Value = $vals[$i]
    }
}

# --- Swap content of rows 23 and 24 ---
Set-RowFV $ws 23 @("Suwon Bluewings", 1, "Daejeon", 3, 2.26, "12/03/2023 16:13", 2.24, "19/03/2023 05:59", 3.36, "12/03/2023 16:13", 3.35, "19/03/2023 05:59", 3.35, "12/03/2023 16:13", 3.52, "19/03/2023 05:59", "https://www.betexplorer.com/football/south-korea/k-league-1/suwon-bluewings-daejeon/jFSQtXfc/")
Set-RowFV $ws 24 @("Daegu", 2, "Jeonbuk", 0, 3.03, "12/03/2023 08:42", 4.22, "19/03/2023 05:59", 3.44, "12/03/2023 08:42", 3.5, "19/03/2023 05:58", 2.43, "12/03/2023 08:42", 1.97, "19/03/2023 05:59", "https://www.betexplorer.com/football/south-korea/k-league-1/daegu-jeonbuk/COTMsDui/")

# --- Swap content of rows 61 and 62 ---
Set-RowFV $ws 61 @("Seoul", 1, "Jeonbuk", 1, 2.64, "29/04/2023 09:42", 2.36, "05/05/2023 06:59", 3.37, "29/04/2023 09:42", 3.37, "05/05/2023 06:59", 2.81, "29/04/2023 09:42", 3.26, "05/05/2023 06:59", "https://www.betexplorer.com/football/south-korea/k-league-1/seoul-jeonbuk/tvYLfELO/")
Set-RowFV $ws 62 @("Daegu", 0, "Ulsan Hyundai", 3, 4.15, "30/04/2023 12:12", 4.06, "05/05/2023 06:51", 3.6, "30/04/2023 12:12", 3.65, "05/05/2023 06:51", 1.92, "30/04/2023 12:12", 1.97, "05/05/2023 06:50", "https://www.betexplorer.com/football/south-korea/k-league-1/daegu-ulsan-hyundai/rqC4JcbP/")

# --- Swap content of rows 130 and 132 (row 131 is untouched) ---
Set-RowFV $ws 130 @("Seoul", 7, "Suwon FC", 2, 1.61, "08/07/2023 12:42", 1.58, "12/07/2023 12:20", 4.3, "08/07/2023 12:42", 4.6, "12/07/2023 12:29", 5.5, "08/07/2023 12:42", 5.54, "12/07/2023 12:29", "https://www.betexplorer.com/football/south-korea/k-league-1/seoul-suwon-fc/6oFysX7O/")
Set-RowFV $ws 132 @("Daejeon", 2, "Jeonbuk", 2, 3.15, "09/07/2023 12:12", 3.42, "12/07/2023 12:24", 3.56, "09/07/2023 12:12", 3.58, "12/07/2023 12:24", 2.31, "09/07/2023 12:12", 2.19, "12/07/2023 12:24", "https://www.betexplorer.com/football/south-korea/k-league-1/daejeon-jeonbuk/f18J5GpP/")

# --- Swap content of rows 159 and 160 ---
Set-RowFV $ws 159 @("Gangwon", 1, "Suwon FC", 2, 2.24, "12/08/2023 13:13", 2.13, "19/08/2023 11:54", 3.6, "12/08/2023 13:13", 3.74, "19/08/2023 11:53", 3.18, "12/08/2023 13:13", 3.42, "19/08/2023 11:59", "https://www.betexplorer.com/football/south-korea/k-league-1/gangwon-suwon-fc/04Ie3KKe/")
Set-RowFV $ws 160 @("Ulsan Hyundai", 1, "Jeonbuk", 0, 1.98, "12/08/2023 12:42", 2.31, "19/08/2023 11:59", 3.77, "12/08/2023 12:42", 3.58, "19/08/2023 11:57", 3.69, "12/08/2023 12:42", 3.17, "19/08/2023 11:59", "https://www.betexplorer.com/football/south-korea/k-league-1/ulsan-hyundai-jeonbuk/6sJi405k/")

# --- Append two new rows (181 and 182) at the bottom, matching formatting of row 180 ---
$ws.Range("A180:V180").Copy() | Out-Null
$ws.Range("A181:V182").PasteSpecial(-4122) | Out-Null
$excel.CutCopyMode = 0

Set-RowFull $ws 181 @(180, "south-korea", "k-league-1", "'2023", 45192.29166666666, "Daejeon", 3, "Suwon Bluewings", 1, 2.04, "17/09/2023 08:42", 1.97, "23/09/2023 06:37", 3.64, "17/09/2023 08:42", 3.84, "23/09/2023 06:37", 3.75, "17/09/2023 08:42", 3.82, "23/09/2023 06:37", "https://www.betexplorer.com/football/south-korea/k-league-1/daejeon-suwon-bluewings/neyHRf3D/")
Set-RowFull $ws 182 @(181, "south-korea", "k-league-1", "'2023", 45192.39583333334, "Jeju Utd", 1, "Seoul", 3, 2.62, "17/09/2023 06:12", 2.98, "23/09/2023 09:29", 3.44, "17/09/2023 06:12", 3.47, "23/09/2023 09:29", 2.73, "17/09/2023 06:12", 2.47, "23/09/2023 09:29", "https://www.betexplorer.com/football/south-korea/k-league-1/jeju-utd-seoul/l0FajvtE/")
